$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on D2:E51 so numeric-looking strings
# (e.g. "1.00", "0.999") are stored as literal text, not converted to numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "66.986.72"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "3.279.69"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "577.97"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").Value = "171.79"
$ws.Range("E6").Value = "  -7.46%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.581"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").Value = "3.278.02"
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  -4.70%  "
$ws.Range("D11").Value = "0.573"
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("D12").Value = "45.14"
$ws.Range("E12").Value = "  -4.33%  "
$ws.Range("D13").Value = "0.0000266"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").Value = "691.16"
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").Value = "3.811.69"
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").Value = "8.28"
$ws.Range("E16").Value = "  -2.67%  "
$ws.Range("D17").Value = "67.111.58"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").Value = "3.289.06"
$ws.Range("E19").Value = "  -1.27%  "
$ws.Range("D20").Value = "17.26"
$ws.Range("E20").Value = "  -3.70%  "
$ws.Range("D21").Value = "10.69"
$ws.Range("E21").Value = "  -3.90%  "
$ws.Range("D22").Value = "0.885"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("D23").Value = "16.86"
$ws.Range("E23").Value = "  -6.11%  "
$ws.Range("E24").Value = "  +3.01%  "
$ws.Range("D25").Value = "99.42"
$ws.Range("E25").Value = "  -3.65%  "
$ws.Range("D26").Value = "3.86"
$ws.Range("E26").Value = "  -2.98%  "
$ws.Range("D27").Value = "2.67"
$ws.Range("E27").Value = "  -4.03%  "
$ws.Range("D28").Value = "33.68"
$ws.Range("E28").Value = "  +3.25%  "
$ws.Range("D29").Value = "9.14"
$ws.Range("E29").Value = "  -4.15%  "
$ws.Range("D30").Value = "8.35"
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("D31").Value = "6.64"
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("D32").Value = "567.83"
$ws.Range("E32").Value = "  -6.80%  "
$ws.Range("D33").Value = "3.846.45"
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("D34").Value = "10.80"
$ws.Range("E34").Value = "  -1.99%  "
$ws.Range("D35").Value = "0.102"
$ws.Range("E35").Value = "  -2.85%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "3.31"
$ws.Range("E37").Value = "  -16.06%  "
$ws.Range("D38").Value = "55.14"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").Value = "0.128"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").Value = "3.44"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("D41").Value = "2.57"
$ws.Range("E41").Value = "  -4.39%  "
$ws.Range("D42").Value = "31.39"
$ws.Range("E42").Value = "  -4.36%  "
$ws.Range("D43").Value = "0.0₃0666"
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "2.99"
$ws.Range("E44").Value = "  -6.55%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.325"
$ws.Range("E45").Value = "  -3.63%  "
$ws.Range("D46").Value = "0.0403"
$ws.Range("E46").Value = "  -2.98%  "
$ws.Range("D47").Value = "0.127"
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").Value = "2.52"
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("D50").Value = "1.37"
$ws.Range("E50").Value = "  +4.54%  "
$ws.Range("D51").Value = "129.84"
$ws.Range("E51").Value = "  -0.88%  "

# Restore original (default/no explicit number format) styling on the data range
$dataRange.ClearFormats()

